$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Date header
Replace-Text "2024-05-04 Saturday" "2024-05-05 Sunday"

# Multiplication problems (old -> new)
Replace-Text "811×8=" "446×3="
Replace-Text "483×9=" "875×6="
Replace-Text "383×9=" "154×6="
Replace-Text "895×9=" "237×8="
Replace-Text "146×6=" "985×2="
Replace-Text "251×8=" "524×8="
Replace-Text "827×4=" "553×5="
Replace-Text "923×5=" "237×9="
Replace-Text "498×5=" "908×6="
Replace-Text "914×8=" "868×4="
Replace-Text "823×2=" "761×2="
Replace-Text "873×5=" "354×6="
Replace-Text "466×8=" "178×2="
Replace-Text "383×7=" "867×4="
Replace-Text "841×8=" "665×7="
Replace-Text "514×4=" "561×7="
Replace-Text "628×6=" "372×8="
Replace-Text "816×2=" "825×3="
Replace-Text "629×4=" "575×5="
Replace-Text "499×9=" "450×7="
Replace-Text "221×9=" "424×3="
Replace-Text "249×3=" "783×9="
Replace-Text "983×7=" "157×9="
Replace-Text "686×5=" "945×6="
Replace-Text "541×3=" "373×5="
